# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  (was "Office Theme")   -> becomes the "Integral" / "Red Violet" theme
#   ppt/theme/theme2.xml  (was "Integral"/"Red Violet", and is the theme actually
#                           applied to the slide master / the whole deck)         -> becomes "Office Theme"
#
# theme2.xml is the one that is reachable through the PowerPoint object model
# (it backs $p.SlideMaster.Theme, which is the only theme/master surface this
# host exposes), so we recolor it to the plain "Office Theme" palette that the
# target XML calls for.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

function RgbVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$cs.Item(1).RGB  = RgbVal 0x00 0x00 0x00   # dk1      000000
$cs.Item(2).RGB  = RgbVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$cs.Item(3).RGB  = RgbVal 0x44 0x54 0x6A   # dk2      44546A
$cs.Item(4).RGB  = RgbVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$cs.Item(5).RGB  = RgbVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$cs.Item(6).RGB  = RgbVal 0xED 0x7D 0x31   # accent2  ED7D31
$cs.Item(7).RGB  = RgbVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$cs.Item(8).RGB  = RgbVal 0xFF 0xC0 0x00   # accent4  FFC000
$cs.Item(9).RGB  = RgbVal 0x44 0x72 0xC4   # accent5  4472C4
$cs.Item(10).RGB = RgbVal 0x70 0xAD 0x47   # accent6  70AD47
$cs.Item(11).RGB = RgbVal 0x05 0x63 0xC1   # hlink    0563C1
$cs.Item(12).RGB = RgbVal 0x95 0x4F 0x72   # folHlink 954F72
